$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ42433925",
    "summ42902547",
    "summ43339800",
    "summ43785870",
    "summ44215172",
    "summ44658666",
    "summ45110798",
    "summ45553137",
    "summ46001949"
)

for ($i = 0; $i -lt $newNames.Length; $i++) {
    $idx = $i + 1
    $ws = $wb.Worksheets.Item($idx)
    $ws.Name = $newNames[$i]
}
